$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly-label entry in column W (adds a new shared string "24-30-aug")
$ws.Range("W30").Value = "24-30-aug"

# Append 4 new daily rows (194-197) covering 2020-08-24 .. 2020-08-27.
# First copy formats from the last existing data row (193) so the new rows
# pick up the same cell styles used for the rest of the table.
$ws.Range("A193:V193").Copy()
$ws.Range("A194:V197").PasteSpecial(-4122)

$rows = @(
  @(44067,192,104472,1809,21983,847,13,5122,36,883,613,0,25495,242,3105,733,1,79330,825,40504,3309,37),
  @(44068,193,106460,1988,21779,859,12,5215,24,885,614,1,25706,211,3085,733,0,80390,1060,41207,3367,58),
  @(44069,194,108403,1943,21062,875,16,5288,73,940,614,0,26033,327,3155,733,0,81646,1256,41939,3421,54),
  @(44070,195,110403,2000,21793,884,9,5379,91,1008,614,0,26361,328,3311,733,0,83150,1504,43014,3459,38)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item(194 + $r, $c + 1).Value = $vals[$c]
    }
}

# Update the view: scroll/selection moves to Y16 (instead of X33 with a pinned
# top-left cell at E1).
$ws.Range("Y16").Select()
